$p = $ppt.ActivePresentation
$s = $p.Slides.Item(15)
$sh = $s.Shapes.Item(2)
$tf = $sh.TextFrame
$tr = $tf.TextRange
$para = $tr.Paragraphs(3)

# Before: 5 runs -
#   1 "Please refer document Parameters "
#   2 "in "
#   3 "Monitor Collections.docx"
#   4 " "
#   5 "for screen prints"
# After: 3 runs -
#   1 "Please refer document Parameters in Monitor Collections.docx"
#   2 " "
#   3 "for screen prints"

# Merge the first three runs' text into run 1.
$para.Runs(1).Text = "Please refer document Parameters in Monitor Collections.docx"

# Remove the now-redundant "in " run (was run 2).
$para.Runs(2).Text = ""

# Remove the now-redundant "Monitor Collections.docx" run (shifted into slot 2).
$para.Runs(2).Text = ""
